$wb = $excel.ActiveWorkbook

# --- 1. Rename the existing "autentification" sheet to "autentification1" ---
$orig = $wb.Worksheets.Item("autentification")
$orig.Name = "autentification1"

# --- 2. Duplicate it (placed right after) to become the new "autentification" sheet ---
$orig.Copy($null, $orig)
$newSheet = $wb.Worksheets.Item($orig.Index + 1)
$newSheet.Name = "autentification"

# --- 3. Trim the new sheet down to just the header + 3 data rows (drop the old blank rows 5:7) ---
$newSheet.Rows("5:7").Delete()

# --- 4. Update the new sheet's header cell A1 to the new "Name" label ---
$newSheet.Range("A1").Value = "Name"

# --- 5. Re-size the new sheet's columns / rows to their final layout ---
$newSheet.Columns.Item(1).ColumnWidth = 39.666666666666664
$newSheet.Columns.Item(2).ColumnWidth = 26.333333333333332
$newSheet.Columns.Item(3).ColumnWidth = 25.833333333333332
$newSheet.Columns.Item(4).ColumnWidth = 28.166666666666664

$newSheet.Rows.Item(2).RowHeight = 38.25
$newSheet.Rows.Item(3).RowHeight = 26.25
$newSheet.Rows.Item(4).RowHeight = 75

# --- 6. Update the selection on the renamed original sheet (no more single active-cell, full data range instead) ---
$orig.Range("A1:D4").Select()

# --- 7. Selection/active cell on the new sheet (left as the active tab) ---
$newSheet.Range("A2").Select()
